$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 16766.666
$ws.Range("I2").Value = 25050
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 25050
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = -24937
$ws.Range("N2").Value = -426

$ws.Range("H12").Value = 640.7692
$ws.Range("I12").Value = 567.44446
$ws.Range("K12").Value = 567.44446
$ws.Range("M12").Value = -397.44446

$ws.Range("H17").Value = 1311
$ws.Range("J17").Value = 1311
$ws.Range("L17").Value = 3933
$ws.Range("N17").Value = -4269

$ws.Range("H18").Value = 1104.25
$ws.Range("I18").Value = 1104.25
$ws.Range("K18").Value = 1104.25
$ws.Range("M18").Value = -820.25

$ws.Range("H19").Value = 399.83334
$ws.Range("I19").Value = 366.33334
$ws.Range("J19").Value = 433.33334
$ws.Range("K19").Value = 366.33334
$ws.Range("L19").Value = 433.33334
$ws.Range("M19").Value = -191.33334
$ws.Range("N19").Value = -783.33334

$ws.Range("H137").Value = 1834.1111
$ws.Range("I137").Value = 807
$ws.Range("K137").Value = 2421
$ws.Range("M137").Value = 129

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1108.6
$ws.Range("I2").Value = 1044.9286
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1044.9286
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -931.9286
$ws.Range("N2").Value = -2226

$ws.Range("H116").Value = 1108.6
$ws.Range("I116").Value = 1044.9286
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1044.9286
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1249.0714
$ws.Range("N116").Value = -6588

$ws.Range("H122").Value = 2017
$ws.Range("I122").Value = 2022.5
$ws.Range("K122").Value = 6067.5
$ws.Range("M122").Value = -3617.5

$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("K132").Value = 10500
$ws.Range("M132").Value = -7970

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1108.6
$ws.Range("I3").Value = 1044.9286
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1044.9286
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -930.9286
$ws.Range("N3").Value = -2228

$ws.Range("H22").Value = 403.8889
$ws.Range("I22").Value = 266.875
$ws.Range("K22").Value = 266.875
$ws.Range("M22").Value = -93.875

$ws.Range("H86").Value = 8200
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 9750
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 9750
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -11996

$ws.Range("H88").Value = 18749.166
$ws.Range("J88").Value = 18749.166
$ws.Range("L88").Value = 18749.166
$ws.Range("N88").Value = -19561.166

$ws.Range("H89").Value = 8200
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 9750
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 48750
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -59982

$ws.Range("H91").Value = 18749.166
$ws.Range("J91").Value = 18749.166
$ws.Range("L91").Value = 18749.166
$ws.Range("N91").Value = -21557.166

$ws.Range("H94").Value = 755.125
$ws.Range("I94").Value = 755.125
$ws.Range("K94").Value = 755.125
$ws.Range("M94").Value = -304.125

$ws.Range("H99").Value = 4010
$ws.Range("I99").Value = 4010
$ws.Range("K99").Value = 4010
$ws.Range("M99").Value = -2512

$ws.Range("H105").Value = 1307
$ws.Range("I105").Value = 1224.6154
$ws.Range("K105").Value = 1224.6154
$ws.Range("M105").Value = 522.3846000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 3999.6667
$ws.Range("I25").Value = 999
$ws.Range("K25").Value = 999
$ws.Range("M25").Value = -825

$ws.Range("H31").Value = 6270.7144
$ws.Range("I31").Value = 2139.5715
$ws.Range("J31").Value = 8336.286
$ws.Range("K31").Value = 2139.5715
$ws.Range("L31").Value = 8336.286
$ws.Range("M31").Value = -1844.5715
$ws.Range("N31").Value = -8926.286

$ws.Range("H34").Value = 6270.7144
$ws.Range("I34").Value = 2139.5715
$ws.Range("J34").Value = 8336.286
$ws.Range("K34").Value = 2139.5715
$ws.Range("L34").Value = 8336.286
$ws.Range("M34").Value = -1937.5715
$ws.Range("N34").Value = -8740.286

$ws.Range("H41").Value = 39473.89
$ws.Range("J41").Value = 45895
$ws.Range("L41").Value = 45895
$ws.Range("N41").Value = -46751

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = $null

$ws.Range("H51").Value = 98765
$ws.Range("J51").Value = 98765
$ws.Range("L51").Value = 98765
$ws.Range("N51").Value = -100237

$ws.Range("H59").Value = 86950
$ws.Range("I59").Value = 86900
$ws.Range("K59").Value = 86900
$ws.Range("M59").Value = -85755

$ws.Range("H60").Value = 57437.145
$ws.Range("J60").Value = 67652
$ws.Range("L60").Value = 67652
$ws.Range("N60").Value = -68674

$ws.Range("H61").Value = 98765
$ws.Range("J61").Value = 98765
$ws.Range("L61").Value = 98765
$ws.Range("N61").Value = -99461

$ws.Range("H86").Value = 3166.6667
$ws.Range("I86").Value = 3250
$ws.Range("K86").Value = 3250
$ws.Range("M86").Value = -2127

$ws.Range("H88").Value = 19999.5
$ws.Range("J88").Value = 19999.5
$ws.Range("L88").Value = 19999.5
$ws.Range("N88").Value = -20811.5

$ws.Range("H89").Value = 3166.6667
$ws.Range("I89").Value = 3250
$ws.Range("K89").Value = 16250
$ws.Range("M89").Value = -10634

$ws.Range("H91").Value = 19999.5
$ws.Range("J91").Value = 19999.5
$ws.Range("L91").Value = 19999.5
$ws.Range("N91").Value = -22807.5

$ws.Range("H95").Value = 24500
$ws.Range("J95").Value = 24500
$ws.Range("L95").Value = 24500
$ws.Range("N95").Value = -29992

$ws.Range("H107").Value = 630.4666999999999
$ws.Range("I107").Value = 310.0909
$ws.Range("J107").Value = 1511.5
$ws.Range("K107").Value = 310.0909
$ws.Range("L107").Value = 1511.5
$ws.Range("M107").Value = 1609.9091
$ws.Range("N107").Value = -5351.5

$ws.Range("H132").Value = 1863.1333
$ws.Range("I132").Value = 1863.1333
$ws.Range("K132").Value = 5589.3999
$ws.Range("M132").Value = -3059.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 456.82608
$ws.Range("J38").Value = 529.8
$ws.Range("L38").Value = 1589.4
$ws.Range("N38").Value = -2283.4

$ws.Range("H57").Value = 2750
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = $null

$ws.Range("H81").Value = 17665
$ws.Range("J81").Value = 17665
$ws.Range("L81").Value = 52995
$ws.Range("N81").Value = -55241

$ws.Range("H84").Value = 17665
$ws.Range("J84").Value = 17665
$ws.Range("L84").Value = 158985
$ws.Range("N84").Value = -170217

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2686.3333
$ws.Range("I126").Value = 2030.5454
$ws.Range("K126").Value = 6091.6362
$ws.Range("M126").Value = -3621.6362

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 300
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 300
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -130
$ws.Range("N16").Value = $null

$ws.Range("H22").Value = 2225.1667
$ws.Range("I22").Value = 2212.75
$ws.Range("K22").Value = 2212.75
$ws.Range("M22").Value = -1917.75

$ws.Range("H26").Value = 700
$ws.Range("I26").Value = 700
$ws.Range("K26").Value = 700
$ws.Range("M26").Value = -405

$ws.Range("H27").Value = 2225.1667
$ws.Range("I27").Value = 2212.75
$ws.Range("K27").Value = 2212.75
$ws.Range("M27").Value = -2105.75

$ws.Range("H100").Value = 5799.8887
$ws.Range("I100").Value = 2439.8
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 2439.8
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -1898.8
$ws.Range("N100").Value = -11082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1700.8334
$ws.Range("I17").Value = 1700.8334
$ws.Range("K17").Value = 1700.8334
$ws.Range("M17").Value = -1528.8334

$ws.Range("H26").Value = 1200
$ws.Range("I26").Value = 1200
$ws.Range("K26").Value = 1200
$ws.Range("M26").Value = -907

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = $null

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = $null

$ws.Range("H81").Value = 7675.5
$ws.Range("J81").Value = 15000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -32122

$ws.Range("H84").Value = 7675.5
$ws.Range("J84").Value = 15000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -160608

$ws.Range("H136").Value = 2966.9092
$ws.Range("I136").Value = 2499
$ws.Range("K136").Value = 7497
$ws.Range("M136").Value = -4947
